# RADSAT-SK Timer BOM update
# "Sourced coin cell battery and tantalum caps, still need to footprint coin cell holders"
#
# This script fills in vendor/pricing info discovered for a few BOM rows:
#   - Row 8  (BT201, BT202 - coin cell battery holder): vendor link + price
#   - Row 9  (C101... - 330uF tantalum caps): corrected footprint/description + vendor link + price
#   - Row 15 (D101 - LED): corrected description + vendor link + price
#   - Row 27 (R102... - 100 ohm resistors): vendor link + price
#
# The "Price total" (M) column is a calculated table column
# (Table1[[#This Row],[Qnty]]*(Table1[[#This Row],[Price]]+Table1[[#This Row],[Price2]]))
# and the "Total Component Cost:" cell (B6) is =SUM(M8:M41), so both recalculate
# automatically once the Price/Price2 cells below are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: BT201, BT202 -- coin cell battery ---
$ws.Range("H8").Value = "https://www.digikey.ca/en/products/detail/keystone-electronics/1061/303558"
$ws.Range("I8").Value = "didkey/mouser"
$ws.Range("J8").Value = 4.09
$ws.Range("K8").Value = "https://www.mouser.ca/ProductDetail/Panasonic-Industrial-Devices/CR-2025-F2N?qs=FNor9lU6pf%2F%2FKr9xpq%252B%252BCQ%3D%3D"
$ws.Range("L8").Value = 1.61

# --- Row 9: C101, C102, C105, C106, C107, C108 -- 330uF tantalum capacitors ---
$ws.Range("E9").Value = "Capacitor_Tantalum_SMD:CP_EIA-7343-30_AVX-N_Pad2.25x2.55mm_HandSolder"
$ws.Range("F9").Value = "330 µF Molded Tantalum Capacitors 10 V 2917 (7343 Metric) 500mOhm @ 100kHz"
$ws.Range("H9").Value = "https://www.digikey.ca/en/products/detail/kyocera-avx/F931A337KNC/4005190"
$ws.Range("I9").Value = "digikey"
$ws.Range("J9").Value = 3.36

# --- Row 15: D101 -- LED ---
$ws.Range("F15").Value = "Green 571nm LED Indication - Discrete 2V 0603 (1608 Metric)"
$ws.Range("H15").Value = "https://www.digikey.ca/en/products/detail/lite-on-inc/LTST-C190KGKT/386815"
$ws.Range("I15").Value = "digikey"
$ws.Range("J15").Value = 0.33

# --- Row 27: R102, R111, R115, R116, R133, R134, R202, R204, R209, R211, R213, R215 -- 100 ohm resistors ---
$ws.Range("H27").Value = "https://www.digikey.ca/en/products/detail/stackpole-electronics-inc/RNCP0805FTD100R/2240209"
$ws.Range("I27").Value = "digikey"
$ws.Range("J27").Value = 0.107

# --- Window / view state (cosmetic, matches author's saved scroll position) ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A10").Select()
